$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.846.12'
$ws.Range("E2").Value = '  -2.99%  '

# Row 3
$ws.Range("D3").Value = '3.426.30'
$ws.Range("E3").Value = '  -4.41%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '185.03'
$ws.Range("E5").Value = '  -8.37%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '534.69'
$ws.Range("E6").Value = '  -4.72%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.614'
$ws.Range("E7").Value = '  -0.72%  '

# Row 8
$ws.Range("D8").Value = '3.418.08'
$ws.Range("E8").Value = '  -4.53%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.09%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.637'
$ws.Range("E10").Value = '  -4.54%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.07'
$ws.Range("E11").Value = '  -1.86%  '

# Row 12
$ws.Range("E12").Value = '  -9.43%  '

# Row 13
$ws.Range("E13").Value = '  -8.35%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.42'
$ws.Range("E14").Value = '  -5.19%  '

# Row 15
$ws.Range("D15").Value = '3.961.62'
$ws.Range("E15").Value = '  -5.26%  '

# Row 16
$ws.Range("E16").Value = '  -2.25%  '

# Row 17
$ws.Range("D17").Value = '3.417.84'
$ws.Range("E17").Value = '  -5.11%  '

# Row 18
$ws.Range("D18").Value = '65.558.17'
$ws.Range("E18").Value = '  -3.19%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.76'
$ws.Range("E19").Value = '  -5.76%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.37'
$ws.Range("E20").Value = '  -7.54%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.988'
$ws.Range("E21").Value = '  -7.54%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '378.92'
$ws.Range("E22").Value = '  -5.24%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '83.09'
$ws.Range("E23").Value = '  -2.16%  '

# Row 24
$ws.Range("E24").Value = '  -8.11%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.05'
$ws.Range("E25").Value = '  -14.46%  '

# Row 26
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.71'
$ws.Range("E26").Value = '  -5.50%  '

# Row 27
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.86'
$ws.Range("E27").Value = '  -5.42%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.69'
$ws.Range("E28").Value = '  -6.82%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.68'
$ws.Range("E29").Value = '  -6.71%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '687.04'
$ws.Range("E30").Value = '  +3.73%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.13'
$ws.Range("E31").Value = '  -3.97%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.91'
$ws.Range("E32").Value = '  -16.58%  '

# Row 33
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.36'
$ws.Range("E33").Value = '  -6.48%  '

# Row 34
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '61.92'
$ws.Range("E34").Value = '  -2.31%  '

# Row 35
$ws.Range("E35").Value = '  -5.28%  '

# Row 36
$ws.Range("E36").Value = '  +0.00%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '37.11'
$ws.Range("E37").Value = '  -11.21%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.391'
$ws.Range("E38").Value = '  -7.21%  '

# Row 39
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.01%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.130'
$ws.Range("E40").Value = '  -4.24%  '

# Row 41
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '29.21'
$ws.Range("E41").Value = '  +31.17%  '

# Row 42
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.907.64'
$ws.Range("E42").Value = '  -11.51%  '

# Row 43
$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.84'
$ws.Range("E43").Value = '  -10.25%  '

# Row 44
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0638'
$ws.Range("E44").Value = '  -15.66%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.67'
$ws.Range("E45").Value = '  -2.48%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0401'
$ws.Range("E46").Value = '  -3.63%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  -12.39%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.128'
$ws.Range("E48").Value = '  -2.21%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '136.93'
$ws.Range("E49").Value = '  -1.27%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.88'
$ws.Range("E50").Value = '  -6.33%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.68'
$ws.Range("E51").Value = '  -1.13%  '
